$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape
# Row 2
$ws.Range("D2").Value = "43.813.68"
$ws.Range("E2").Value = "  -0.70%  "
# Row 3
$ws.Range("D3").Value = "2.232.89"
$ws.Range("E3").Value = "  -2.45%  "
# Row 4
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
$ws.Range("D5").Value = "'0.647"
$ws.Range("E5").Value = "  +4.16%  "
# Row 6
$ws.Range("D6").Value = "'229.80"
$ws.Range("E6").Value = "  -0.54%  "
# Row 7
$ws.Range("D7").Value = "'62.63"
$ws.Range("E7").Value = "  +2.79%  "
# Row 8
$ws.Range("E8").Value = "  +0.03%  "
# Row 9
$ws.Range("D9").Value = "'0.446"
$ws.Range("E9").Value = "  +4.78%  "
# Row 10
$ws.Range("D10").Value = "'0.0956"
$ws.Range("E10").Value = "  +1.63%  "
# Row 11
$ws.Range("D11").Value = "'56.80"
$ws.Range("E11").Value = "  -1.88%  "
# Row 12
$ws.Range("D12").Value = "'26.28"
$ws.Range("E12").Value = "  +8.28%  "
# Row 13
$ws.Range("E13").Value = "  +1.07%  "
# Row 14
$ws.Range("D14").Value = "2.564.74"
$ws.Range("E14").Value = "  -2.51%  "
# Row 15
$ws.Range("D15").Value = "'15.35"
$ws.Range("E15").Value = "  -2.49%  "
# Row 16
$ws.Range("D16").Value = "'6.08"
$ws.Range("E16").Value = "  +2.99%  "
# Row 17
$ws.Range("D17").Value = "'0.819"
$ws.Range("E17").Value = "  +0.68%  "
# Row 18
$ws.Range("D18").Value = "2.228.05"
$ws.Range("E18").Value = "  -2.49%  "
# Row 19
$ws.Range("D19").Value = "43.632.19"
$ws.Range("E19").Value = "  -0.71%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  +3.38%  "
# Row 21
$ws.Range("D21").Value = "'72.52"
$ws.Range("E21").Value = "  -1.61%  "
# Row 22
$ws.Range("D22").Value = "'5.99"
$ws.Range("E22").Value = "  -4.09%  "
# Row 23
$ws.Range("D23").Value = "'246.99"
$ws.Range("E23").Value = "  -2.56%  "
# Row 24
$ws.Range("E24").Value = "  -0.02%  "
# Row 25
$ws.Range("E25").Value = "  -6.14%  "
# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.29"
$ws.Range("E26").Value = "  -2.96%  "
# Row 27
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "'3.37"
$ws.Range("E27").Value = "  +22.66%  "
# Row 28
$ws.Range("D28").Value = "'9.79"
$ws.Range("E28").Value = "  -0.99%  "
# Row 29
$ws.Range("D29").Value = "'170.28"
$ws.Range("E29").Value = "  -0.48%  "
# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = "  -1.71%  "
# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'20.61"
$ws.Range("E31").Value = "  +0.22%  "
# Row 32
$ws.Range("E32").Value = "  -4.05%  "
# Row 33
$ws.Range("E33").Value = "  +3.59%  "
# Row 34
$ws.Range("D34").Value = "'0.0689"
$ws.Range("E34").Value = "  +4.60%  "
# Row 35
$ws.Range("D35").Value = "'4.73"
$ws.Range("E35").Value = "  -1.36%  "
# Row 36
$ws.Range("D36").Value = "'4.87"
$ws.Range("E36").Value = "  -3.51%  "
# Row 37
$ws.Range("D37").Value = "'3.63"
$ws.Range("E37").Value = "  -0.57%  "
# Row 38
$ws.Range("E38").Value = "  -2.40%  "
# Row 39
$ws.Range("D39").Value = "'2.26"
$ws.Range("E39").Value = "  -5.83%  "
# Row 40
$ws.Range("D40").Value = "'0.0256"
$ws.Range("E40").Value = "  +2.37%  "
# Row 41
$ws.Range("E41").Value = "  -0.02%  "
# Row 42
$ws.Range("D42").Value = "'0.000217"
$ws.Range("E42").Value = "  -2.26%  "
# Row 43
$ws.Range("D43").Value = "'8.20"
$ws.Range("E43").Value = "  -6.10%  "
# Row 44
$ws.Range("D44").Value = "'17.00"
$ws.Range("E44").Value = "  -0.31%  "
# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0945"
$ws.Range("E45").Value = "  -2.75%  "
# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'96.41"
$ws.Range("E46").Value = "  -2.36%  "
# Row 47
$ws.Range("E47").Value = "  -2.45%  "
# Row 48
$ws.Range("D48").Value = "'4.33"
$ws.Range("E48").Value = "  -1.49%  "
# Row 49
$ws.Range("E49").Value = "  +2.01%  "
# Row 50
$ws.Range("D50").Value = "1.425.67"
$ws.Range("E50").Value = "  -3.59%  "
# Row 51
$ws.Range("E51").Value = "  +1.51%  "
